## Apply the edits described by the commit:
##  - Supplier Y's second invoice number (rows 8-9) changes from 50-9501 to 50-9505
##  - Column D (Cost) gets a custom currency number format
##  - Column D is auto-fit to its new (wider) content
##  - Selection moves to A10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supplier Y's invoice number for the 6650 / 125 rows (8 and 9) changes
$ws.Range("B8").Value = "50-9505"
$ws.Range("B9").Value = "50-9505"

# Apply a custom accounting/currency number format to the Cost column (D2:D13)
$costFormat = '_-[$$-409]* #,##0.00_ ;_-[$$-409]* \-#,##0.00\ ;_-[$$-409]* "-"??_ ;_-@_ '
$ws.Range("D2:D13").NumberFormat = $costFormat

# Column D is now wider to fit the new currency-formatted values
$ws.Columns("D").AutoFit()

# Move the active selection to A10
[void]$ws.Range("A10").Select()
